# FUNCTIONALITY: Wrote two new test cases.
# Update the "ListView" test-suite row (row 9) to account for two newly
# written test cases: one more automated test case and one more total
# test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9 ("ListView"): Automated Test Cases (B9) and Total Test Cases (C9)
# both increase by one, from 4 to 5.
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 5

# Move the active selection from E9 to D9, matching the saved sheet view.
$ws.Range("D9").Select()
